$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need an explicit
# text number-format so Excel does not silently convert them to numeric values.
$ws.Range("D2").Value = '43.487.70'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '2.335.72'
$ws.Range("E3").Value = '  -1.68%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.32'
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.53'
$ws.Range("E6").Value = '  -3.38%  '
$ws.Range("E7").Value = '  -1.76%  '
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.19'
$ws.Range("E10").Value = '  -2.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.59'
$ws.Range("E11").Value = '  -3.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0797'
$ws.Range("E12").Value = '  -1.92%  '
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("E14").Value = '  -2.70%  '
$ws.Range("D15").Value = '2.697.12'
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.71'
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").Value = '2.327.15'
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").Value = '43.392.61'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.78'
$ws.Range("E20").Value = '  -1.78%  '
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("E22").Value = '  -2.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.88'
$ws.Range("E23").Value = '  -0.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '238.48'
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("E25").Value = '  -3.69%  '
$ws.Range("E26").Value = '  -3.63%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.02'
$ws.Range("E28").Value = '  -3.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.62'
$ws.Range("E29").Value = '  -6.67%  '
$ws.Range("E30").Value = '  -2.03%  '
$ws.Range("E31").Value = '  -3.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.72'
$ws.Range("E32").Value = '  +1.80%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.07'
$ws.Range("E34").Value = '  -3.74%  '
$ws.Range("E35").Value = '  -5.04%  '
$ws.Range("E36").Value = '  -6.51%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0709'
$ws.Range("E37").Value = '  -4.75%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '16.86'
$ws.Range("E38").Value = '  -7.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.91'
$ws.Range("E39").Value = '  -6.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.83'
$ws.Range("E40").Value = '  -6.30%  '
$ws.Range("E41").Value = '  -3.11%  '
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.41'
$ws.Range("E43").Value = '  -3.24%  '
$ws.Range("D44").Value = '1.979.90'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("E45").Value = '  -1.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.64'
$ws.Range("E46").Value = '  -7.26%  '
$ws.Range("E47").Value = '  -6.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.91'
$ws.Range("E48").Value = '  -5.82%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.91'
$ws.Range("E49").Value = '  +4.09%  '
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.65'
$ws.Range("E50").Value = '  -4.36%  '
$ws.Range("D51").Value = '2.561.46'
$ws.Range("E51").Value = '  +0.23%  '
